$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (existing entry's data changed) ---
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 43734
$ws.Range("C2").Value = "no"
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = "looking at phone"
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 5

# --- Update row 3 (existing entry's data changed) ---
$ws.Range("A3").Value = 0
$ws.Range("F3").Value = 15
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3

# --- Add new rows 4-7, cloning the formatting of row 3 first ---
$ws.Range("A3:O3").Copy() | Out-Null
$ws.Range("A4:O4").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:O5").PasteSpecial(-4122) | Out-Null
$ws.Range("A6:O6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:O7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 43732
$ws.Range("C4").Value = "yes"
$ws.Range("D4").Value = "no"
$ws.Range("E4").Value = "no"
$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = "meditation"
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 4

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 43732
$ws.Range("C5").Value = "yes"
$ws.Range("D5").Value = "no"
$ws.Range("E5").Value = "no"
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = "exercise"
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 4

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 43732
$ws.Range("C6").Value = "yes"
$ws.Range("D6").Value = "no"
$ws.Range("E6").Value = "no"
$ws.Range("F6").Value = 45
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = "looking at phone"
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 2
$ws.Range("N6").Value = 4
$ws.Range("O6").Value = 3

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 43732
$ws.Range("C7").Value = "yes"
$ws.Range("D7").Value = "no"
$ws.Range("E7").Value = "no"
$ws.Range("F7").Value = 30
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "exercise"
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 3
$ws.Range("O7").Value = 2
